# Generate Report for Handback
#
# This mirrors the localization-status report-generation step: the status
# text flips from "Ready for handoff" to "Handed back: in sync with en-US"
# on every sheet, and each per-locale sheet gets its "Latest Target File",
# "Latest Handback File" and "Latest Handback DateTime" columns filled in
# (plus a hyperlink on the newly-populated target-file cell). Column widths
# on the status/filename columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFile     = "88397580-d60d-44c7-8f0a-16fb571530b6.md"
$mdDisplay  = "88397580-d60d-44c7-8f0a-16fb571530b6.md"

# --- Overview sheet: widen the per-locale status columns and refresh the text
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("I2").Value = $mdFile
$zh.Range("J2").Value = "88397580-d60d-44c7-8f0a-16fb571530b6.62cded27da37696ec18b99c658c58f586eccb7f4.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-21 15:05:44"
$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a61cf28293c9e30d7c9bb9a49c15384eaa3ed3a/e2e/88397580-d60d-44c7-8f0a-16fb571530b6.md", $null, $null, $mdDisplay) | Out-Null
$zh.Range("I2").Style = "HyperLink"

# --- de-de sheet
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("I2").Value = $mdFile
$de.Range("J2").Value = "88397580-d60d-44c7-8f0a-16fb571530b6.62cded27da37696ec18b99c658c58f586eccb7f4.de-de.xlf"
$de.Range("K2").Value = "2016-08-21 15:05:50"
$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a61cf28293c9e30d7c9bb9a49c15384eaa3ed3a/e2e/88397580-d60d-44c7-8f0a-16fb571530b6.md", $null, $null, $mdDisplay) | Out-Null
$de.Range("I2").Style = "HyperLink"
